$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content/label updates -------------------------------------------------
# (Set in the same order the new shared strings appear so the rebuilt
#  sharedStrings table lines up with the target file.)

# "location" first appears in H16, reused in L16
$ws.Cells.Item(16, 8).Value = "location"        # H16
$ws.Cells.Item(16, 12).Value = "location"       # L16

# "(int,int)" first appears in I16, reused in M16/M17
$ws.Cells.Item(16, 9).Value = "(int,int)"        # I16
$ws.Cells.Item(16, 13).Value = "(int,int)"       # M16
$ws.Cells.Item(17, 13).Value = "(int,int)"       # M17

# new Walker-class method row
$ws.Cells.Item(20, 2).Value = "one_step_random_direction_Walker"   # B20

# new class block title
$ws.Cells.Item(13, 12).Value = "Magic_Gate Class"  # L13

# more new method rows
$ws.Cells.Item(21, 2).Value = "random_step_random_direction"  # B21
$ws.Cells.Item(22, 2).Value = "one_step_discrete_direction"   # B22

# Magic_Gate attribute rows
$ws.Cells.Item(17, 12).Value = "destination"   # L17
$ws.Cells.Item(15, 12).Value = "properties"    # L15
$ws.Cells.Item(15, 13).Value = "Type"          # M15
$ws.Cells.Item(15, 14).Value = "Encapsulation" # N15
$ws.Cells.Item(16, 14).Value = "private"       # N16
$ws.Cells.Item(17, 14).Value = "private"       # N17

# Walker attributes note
$ws.Cells.Item(11, 2).Value = "direction- one of 8 directions"  # B11

# Walker-class-methods notes (H column)
$ws.Cells.Item(17, 8).Value = "length"               # H17
$ws.Cells.Item(18, 8).Value = "check intersections"  # H18

# Attribute renames x / y / z
$ws.Cells.Item(8, 2).Value = "x"                 # B8 (was "direction")
$ws.Cells.Item(9, 2).Value = "y"                 # B9 (was "step_size")
$ws.Cells.Item(10, 2).Value = "z"                # B10 (new row)
$ws.Cells.Item(10, 3).Value = "defaults to 0"    # C10 (new row)

# --- Styling for the new/changed cells ------------------------------------------
$ws.Cells.Item(13, 12).Font.Bold = $true
$ws.Cells.Item(13, 12).Font.Underline = $true

$ws.Cells.Item(15, 12).Font.Bold = $true
$ws.Cells.Item(15, 13).Font.Bold = $true
$ws.Cells.Item(15, 14).Font.Bold = $true

$ws.Cells.Item(20, 2).Font.Bold = $true
$ws.Cells.Item(21, 2).Font.Bold = $true
$ws.Cells.Item(22, 2).Font.Bold = $true

# --- Sheet view: zoom + selection -------------------------------------------------
$excel.ActiveWindow.Zoom = 70
[void]$ws.Range("B37").Select()

# --- Column widths -----------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 31.3
$ws.Columns.Item(8).ColumnWidth = 15.5
